$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.232.16"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").Value = "'1.580.76"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'212.72"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  +6.93%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'26.32"
$ws.Range("E8").Value = "  +10.68%  "
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'1.806.49"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "'1.585.81"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("D14").Value = "'29.251.70"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").Value = "'3.72"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "'0.524"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "'62.82"
$ws.Range("D18").Value = "'237.82"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").Value = "'154.14"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("D27").Value = "'15.15"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'0.0470"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'1.07"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "'1.425.55"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'1.51"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("E37").Value = "  +6.70%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "'1.97"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("D42").Value = "'54.60"
$ws.Range("E42").Value = "  +26.99%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'0.791"
$ws.Range("E44").Value = "  +2.00%  "
$ws.Range("D45").Value = "'0.0472"
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("D46").Value = "'64.62"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").Value = "'5.34"
$ws.Range("D48").Value = "'1.718.47"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'0.840"
$ws.Range("E49").Value = "  -6.73%  "
$ws.Range("D50").Value = "'85.43"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  -2.35%  "
